# Apply revised dose/SD values from the updated poster data for the two
# phase III trials (NCT04072354 and NCT04092686), rows 3-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (NCT04072354, 50mg)
$ws.Range("E3").Value = 19.07
$ws.Range("F3").Value = 142
$ws.Range("H3").Value = 18.06
$ws.Range("K3").Value = 0.127629525937716
$ws.Range("L3").Value = 0.02096920151256492

# Row 4 (NCT04072354, 75mg)
$ws.Range("E4").Value = 19.27
$ws.Range("H4").Value = 18.06
$ws.Range("K4").Value = -0.01583516682992753
$ws.Range("L4").Value = 0.0207860183837838

# Row 5 (NCT04092686, 75mg)
$ws.Range("E5").Value = 18.55
$ws.Range("F5").Value = 153
$ws.Range("H5").Value = 18.67
$ws.Range("K5").Value = -0.1125913698497899
$ws.Range("L5").Value = 0.01955051899631097

# Row 6 (NCT04092686, 100mg)
$ws.Range("E6").Value = 18.49
$ws.Range("F6").Value = 152
$ws.Range("H6").Value = 18.67
$ws.Range("K6").Value = -0.2041690083470224
$ws.Range("L6").Value = 0.01965697560427554
